# Added common test data for Belgium and Germany market
#
# Inserts two new repeater models (P32AR, P32DR) into the "Repeaters" sheet,
# ahead of the existing PR1DS/PR8AS/ZXF/ZXFEV/Wg/Repeaters rows (mirroring
# the data already present on "Repeaters_Updated"), and updates the active
# sheet/selection state to reflect the edit.

$wb = $excel.ActiveWorkbook

$wsRepeaters        = $wb.Worksheets.Item("Repeaters")
$wsRepeatersUpdated = $wb.Worksheets.Item("Repeaters_Updated")

# Push the existing tail (PR1DS, PR8AS, ZXF, ZXFEV, Wg, Repeaters) down by
# two rows, carrying their formatting with them, then reuse the row-15
# formatting for the two freshly inserted rows before writing their values.
$wsRepeaters.Range("A16:A21").Copy($wsRepeaters.Range("A18:A23"))
$wsRepeaters.Range("A15").Copy($wsRepeaters.Range("A16:A17"))
$wsRepeaters.Range("A16").Value = "P32AR"
$wsRepeaters.Range("A17").Value = "P32DR"

# Reflect the edit in the UI state: "Repeaters" becomes the active sheet/tab,
# with A8:A23 selected (anchored at A8); "Repeaters_Updated" keeps its own
# A8:A23 selection without stealing tab focus.
$wsRepeatersUpdated.Range("A8:A23").Select()
$wsRepeaters.Activate()
$wsRepeaters.Range("A8:A23").Select()
